$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 228
$ws.Range("D228").Value = 45041
$ws.Range("K228").Value = "Crimpson Seedless"
$ws.Range("M228").Value = 100
$ws.Range("P228").Value = 9500
$ws.Range("R228").Value = "Región de O'Higgins"
$ws.Range("S228").Value = 528

# Row 229
$ws.Range("D229").Value = 45041
$ws.Range("M229").Value = 150
$ws.Range("P229").Value = 9667
$ws.Range("S229").Value = 537

# Row 230
$ws.Range("D230").Value = 44705
$ws.Range("M230").Value = 220
$ws.Range("P230").Value = 9545
$ws.Range("R230").Value = "Provincia de Limarí"
$ws.Range("S230").Value = 530

# Row 231
$ws.Range("D231").Value = 45015
$ws.Range("L231").Value = "Primera"
$ws.Range("M231").Value = 100
$ws.Range("N231").Value = 9000
$ws.Range("O231").Value = 10000
$ws.Range("P231").Value = 9500
$ws.Range("S231").Value = 528

# Row 232
$ws.Range("D232").Value = 44342
$ws.Range("K232").Value = "Red Globe"
$ws.Range("M232").Value = 100
$ws.Range("N232").Value = 9000
$ws.Range("O232").Value = 10000
$ws.Range("P232").Value = 9500
$ws.Range("Q232").Value = "`$/bandeja 18 kilos"
$ws.Range("R232").Value = "Región de O'Higgins"
$ws.Range("S232").Value = 528
$ws.Range("T232").Value = 18

# Row 233
$ws.Range("D233").Value = 44342
$ws.Range("L233").Value = "Segunda"
$ws.Range("M233").Value = 50
$ws.Range("N233").Value = 8000
$ws.Range("O233").Value = 8000
$ws.Range("P233").Value = 8000
$ws.Range("S233").Value = 444

# Row 234
$ws.Range("A234").Value = 11
$ws.Range("B234").Value = "Vega Monumental Concepción"
$ws.Range("C234").Value = "Bíobío"
$ws.Range("D234").Value = 44217
$ws.Range("D234").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E234").Value = 8
$ws.Range("F234").Value = "Fruta"
$ws.Range("G234").Value = 100109
$ws.Range("H234").Value = "Uva"
$ws.Range("I234").Value = 100109001
$ws.Range("J234").Value = "Uva"
$ws.Range("K234").Value = "Superior Seedless"
$ws.Range("L234").Value = "Primera"
$ws.Range("M234").Value = 200
$ws.Range("N234").Value = 10000
$ws.Range("O234").Value = 11000
$ws.Range("P234").Value = 10500
$ws.Range("Q234").Value = "`$/bandeja 10 kilos"
$ws.Range("R234").Value = "Provincia del Elquí"
$ws.Range("S234").Value = 1050
$ws.Range("T234").Value = 10

# Row 235
$ws.Range("A235").Value = 11
$ws.Range("B235").Value = "Vega Monumental Concepción"
$ws.Range("C235").Value = "Bíobío"
$ws.Range("D235").Value = 45007
$ws.Range("D235").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E235").Value = 8
$ws.Range("F235").Value = "Fruta"
$ws.Range("G235").Value = 100109
$ws.Range("H235").Value = "Uva"
$ws.Range("I235").Value = 100109001
$ws.Range("J235").Value = "Uva"
$ws.Range("K235").Value = "Red Globe"
$ws.Range("L235").Value = "Primera"
$ws.Range("M235").Value = 180
$ws.Range("N235").Value = 10000
$ws.Range("O235").Value = 11000
$ws.Range("P235").Value = 10444
$ws.Range("Q235").Value = "`$/bandeja 18 kilos"
$ws.Range("R235").Value = "Región de O'Higgins"
$ws.Range("S235").Value = 580
$ws.Range("T235").Value = 18
